# Regenerate all penyata to follow new data and format
# --------------------------------------------------------------
# This script updates the HOMEROOM 1MURNI-2023 statement sheet:
#  - relabels the "Kali Pertama/Kedua/Ketiga/Keempat" checkpoints to
#    "Semakan Kali ..." (the label is reused across three review
#    sections, so all 12 occurrences are updated)
#  - normalises the competition-name list to title case and appends
#    three new competitions that were missing
#  - reflows the header ("STATEMENT OF HOMEROOM ACCOUNT" moves from
#    E4 to D4) and updates merged cell ranges to match the new layout
#  - repositions/resizes the letterhead picture
#  - refreshes page-setup/print options used for the regenerated report
# --------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------- 1. Header reflow (row 4 / row 5) ----------------------------
# "STATEMENT OF HOMEROOM ACCOUNT" moves from E4 to D4, and D4:G4 becomes
# the merged title cell.
$titleText = $ws.Range("E4").Value()
$ws.Range("E4").Value = ""
$ws.Range("D4").Value = $titleText
$ws.Range("D4:G4").Merge()

# B5 ("MRSM GEMENCHEH, NEGERI SEMBILAN.") now spans B5:C5.
$ws.Range("B5:C5").Merge()

# ---------- 2. Section heading merge (row 12) ----------------------------
$ws.Range("B12:F12").Merge()

# ---------- 3. Fix up the "Kali ..." merges / totals-row merges ----------
# B15:C15 is no longer merged.
$ws.Range("B15:C15").UnMerge()

# B21:C21, B27:C27, B33:C33 grow to span through column E.
$ws.Range("B21:C21").UnMerge()
$ws.Range("B21:E21").Merge()

$ws.Range("B27:C27").UnMerge()
$ws.Range("B27:E27").Merge()

$ws.Range("B33:C33").UnMerge()
$ws.Range("B33:E33").Merge()

# Grand-total row (43) becomes merged B43:E43.
$ws.Range("B43:E43").Merge()

# ---------- 4. Relabel the three "Kali Pertama/Kedua/Ketiga/Keempat" ------
#              checklists -> "Semakan Kali ..."
$kaliLabels = @("Semakan Kali Pertama", "Semakan Kali Kedua", "Semakan Kali Ketiga", "Semakan Kali Keempat")
$kaliRowGroups = @(16, 22, 28)
foreach ($startRow in $kaliRowGroups) {
    for ($i = 0; $i -lt 4; $i++) {
        $row = $startRow + $i
        $ws.Range("C$row").Value = $kaliLabels[$i]
    }
}

# ---------- 5. Normalise competition names + add the missing ones -------
$ws.Range("C34").Value = "Bouquet Kreatif"
$ws.Range("C35").Value = "Kad Raya Untuk Guruku"
$ws.Range("C36").Value = "Riang Ria Kuih Raya"
$ws.Range("C37").Value = "Creative Collage"
$ws.Range("C38").Value = "Lompat Getah"
$ws.Range("C39").Value = "Theme Party"
$ws.Range("C40").Value = "Hari Koperasi"

# ---------- 6. Letterhead picture: reposition + resize -------------------
$shp = $ws.Shapes.Item(1)
$shp.Left = 41.2125
$shp.Top = 14.25
$shp.Width = 46.5
$shp.Height = 47.25

# ---------- 7. Page setup / print options for the regenerated report -----
$ps = $ws.PageSetup
$ps.FitToPagesWide = 1
$ps.FitToPagesTall = 1
$ps.Zoom = $false
$ps.CenterHorizontally = $true
$ps.HeaderMargin = 0
$ps.FooterMargin = 0
